$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (borders, styles) of the previous data row down into
# the new row so the new cells line up with the existing table (same thin
# borders used by row 2), then overwrite with the new reader's data.
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 44492
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("C3").Value = "Артуром Конаном Дойлом"
$ws.Range("D3").Value = "Шерлок холмс"

[void]$ws.Range("A3:D3").Select()
